# Generate Report for Handback
#
# Rows for "38222f7e-5091-415f-af85-f72e2e08270e.md" (row 3) and its
# dependent "bf9a8cfe-181b-4734-b5d9-426fbaf76a41.md" (row 4) have now
# been handed back, so:
#   - Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview sheet + per-language
#     sheets).
#   - The per-language sheets gain "Latest Target File" / "Latest
#     Handback File" hyperlinked entries (mirroring the source file /
#     handoff file), and the "Latest Handback DateTime" is stamped.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 3 & 4, columns B (zh-cn) and C (de-de)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $handedBack
$wsOverview.Range("C3").Value = $handedBack
$wsOverview.Range("B4").Value = $handedBack
$wsOverview.Range("C4").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 3 - 38222f7e-5091-415f-af85-f72e2e08270e.md
$wsZh.Range("B3").Value = $handedBack
$wsZh.Range("E3").Value = "38222f7e-5091-415f-af85-f72e2e08270e.md"
$wsZh.Range("F3").Value = "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.zh-cn.xlf"
$wsZh.Range("G3").Value = "2016-01-29 02:32:22"

# Row 4 - bf9a8cfe-181b-4734-b5d9-426fbaf76a41.md (dependent on row 3)
$wsZh.Range("B4").Value = $handedBack
$wsZh.Range("E4").Value = "38222f7e-5091-415f-af85-f72e2e08270e.md"
$wsZh.Range("F4").Value = "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.zh-cn.xlf"
$wsZh.Range("G4").Value = "2016-01-29 02:32:22"

$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fd4bec291a8efa2cb7bf145197318302d1024807/e2e/38222f7e-5091-415f-af85-f72e2e08270e.md", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/41a41de4ffc4b03fa16edea015faa4c31849ad33/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.zh-cn.xlf", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fd4bec291a8efa2cb7bf145197318302d1024807/e2e/38222f7e-5091-415f-af85-f72e2e08270e.md", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/41a41de4ffc4b03fa16edea015faa4c31849ad33/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.zh-cn.xlf", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 3 - 38222f7e-5091-415f-af85-f72e2e08270e.md
$wsDe.Range("B3").Value = $handedBack
$wsDe.Range("E3").Value = "38222f7e-5091-415f-af85-f72e2e08270e.md"
$wsDe.Range("F3").Value = "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.de-de.xlf"
$wsDe.Range("G3").Value = "2016-01-29 02:32:41"

# Row 4 - bf9a8cfe-181b-4734-b5d9-426fbaf76a41.md (dependent on row 3)
$wsDe.Range("B4").Value = $handedBack
$wsDe.Range("E4").Value = "38222f7e-5091-415f-af85-f72e2e08270e.md"
$wsDe.Range("F4").Value = "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.de-de.xlf"
$wsDe.Range("G4").Value = "2016-01-29 02:32:41"

$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7de58900d0b4aab5146e04a3a833dd79ef3b6195/e2e/38222f7e-5091-415f-af85-f72e2e08270e.md", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/871ed3c91647206506b5c39a8b961c70c35611c9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/tianzh/38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.de-de.xlf", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7de58900d0b4aab5146e04a3a833dd79ef3b6195/e2e/38222f7e-5091-415f-af85-f72e2e08270e.md", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/871ed3c91647206506b5c39a8b961c70c35611c9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/tianzh/38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.de-de.xlf", "", "", "38222f7e-5091-415f-af85-f72e2e08270e.a21582f05efb00c6c84aa42a4a36c6169c86c91a.de-de.xlf") | Out-Null
